# Add data for 2023-02-21
# Updates 2023 year-to-date violent-crime counts (column J) across the
# "Citywide Totals", "By Neighborhood" and individual per-neighborhood
# sheets. A couple of 2022 (column I) and 2020 (column G) values also
# receive small retroactive corrections, matching the upstream diff.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 872
$ws.Range("J3").Value = 966
$ws.Range("G4").Value = 1454
$ws.Range("I4").Value = 1753
$ws.Range("J4").Value = 214
$ws.Range("J5").Value = 71
$ws.Range("J6").Value = 1360
$ws.Range("G7").Value = 24679
$ws.Range("I7").Value = 26185
$ws.Range("J7").Value = 3483

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 98
$ws.Range("J8").Value = 225
$ws.Range("J17").Value = 6
$ws.Range("J19").Value = 113
$ws.Range("J20").Value = 78
$ws.Range("J29").Value = 179
$ws.Range("J30").Value = 15
$ws.Range("J33").Value = 149
$ws.Range("J34").Value = 23
$ws.Range("J37").Value = 119
$ws.Range("J38").Value = 2
$ws.Range("J41").Value = 24
$ws.Range("G42").Value = 999
$ws.Range("J42").Value = 157
$ws.Range("J43").Value = 41
$ws.Range("J46").Value = 10
$ws.Range("J50").Value = 17
$ws.Range("J52").Value = 77
$ws.Range("I63").Value = 177
$ws.Range("J63").Value = 18
$ws.Range("J64").Value = 20
$ws.Range("J65").Value = 92
$ws.Range("J67").Value = 132
$ws.Range("J72").Value = 13
$ws.Range("J74").Value = 6
$ws.Range("J76").Value = 56
$ws.Range("J77").Value = 31
$ws.Range("J78").Value = 46
$ws.Range("J79").Value = 110
$ws.Range("J83").Value = 79
$ws.Range("J85").Value = 143
$ws.Range("J86").Value = 14
$ws.Range("J89").Value = 40
$ws.Range("J91").Value = 46
$ws.Range("J96").Value = 41
$ws.Range("J97").Value = 20
$ws.Range("G101").Value = 24679
$ws.Range("I101").Value = 26185
$ws.Range("J101").Value = 3483

# Sheet 3: South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 46
$ws.Range("J4").Value = 12
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 143

# Sheet 5: Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 77

# Sheet 7: Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 68
$ws.Range("J3").Value = 75
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 225

# Sheet 9: Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 98

# Sheet 10: Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 40

# Sheet 11: West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 14
$ws.Range("J4").Value = 4
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 41

# Sheet 13: Fuller Park
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J3").Value = 6
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 15

# Sheet 14: Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 32
$ws.Range("J3").Value = 39
$ws.Range("J7").Value = 119

# Sheet 16: North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 56
$ws.Range("J7").Value = 132

# Sheet 19: New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 22
$ws.Range("J7").Value = 92

# Sheet 20: South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 22
$ws.Range("J7").Value = 79

# Sheet 22: Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 149

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 59
$ws.Range("J4").Value = 5
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 179

# Sheet 26: Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 113

# Sheet 29: River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 56

# Sheet 31: Hermosa
$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 24

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 26
$ws.Range("G4").Value = 45
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 97
$ws.Range("G7").Value = 999
$ws.Range("J7").Value = 157

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 46

# Sheet 38: Jefferson Park
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 10

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 46

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 36
$ws.Range("J7").Value = 110

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 20

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J3").Value = 26
$ws.Range("J4").Value = 7
$ws.Range("J5").Value = 1
$ws.Range("J7").Value = 78

# Sheet 46: Burnside
$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("J6").Value = 2
$ws.Range("J7").Value = 6

# Sheet 50: Garfield Ridge
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 23

# Sheet 56: Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 17

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 29

# Sheet 65: West Town
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J2").Value = 5
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 20

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 14

# Sheet 79: Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 41

# Sheet 82: Old Town
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 13

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 31

# Sheet 95: Printers Row
$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 6

# Sheet 100: Grant Park
$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("J4").Value = 1
$ws.Range("J6").Value = 2
